# Ticket 18: replace "2." with "18" + "." + " Исходные положения и основные
# понятия теории познания И. Канта" as three separate runs.
#
# Using TrackRevisions while editing (and then accepting all revisions)
# prevents the engine from silently re-merging adjacent runs that share
# identical (empty) formatting, which is what happens with plain edits.

$d = $word.ActiveDocument
$d.TrackRevisions = $true

# Replace the whole "2." run's text with "18".
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.Text = "18"

# Append "." as its own run right after "18".
$p2 = $d.Paragraphs(1)
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertAfter(".")

# Append the remaining title text as its own run.
$p3 = $d.Paragraphs(1)
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertAfter(" Исходные положения и основные понятия теории познания И. Канта")

$d.TrackRevisions = $false
$d.AcceptAllRevisions()
